# Processed Results - mean and median increase
# Adds "Mean increase" / "Median increase" header + formula cells to each
# of the four frequency groups (Low, Medium, High, All) on Blad1, matching
# the "Mean increase"/"Median increase" columns added in D/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

function Add-IncreaseBlock($HeaderRow, $DataRow, $MeanRef, $MedianRef) {

    # Headers (bold, matching the existing section-header style).
    $hdrMean = $ws.Range("D" + $HeaderRow)
    $hdrMean.Value = "Mean increase"
    $hdrMean.Font.Bold = $true

    $hdrMedian = $ws.Range("F" + $HeaderRow)
    $hdrMedian.Value = "Median increase"
    $hdrMedian.Font.Bold = $true

    # Formula cells - percentage increase relative to a fixed baseline.
    $meanCell = $ws.Range("D" + $DataRow)
    $meanCell.Formula = "= ((" + $MeanRef + " / 95.321842) * 100) - 100"
    $meanCell.ClearFormats()

    $medianCell = $ws.Range("F" + $DataRow)
    $medianCell.Formula = "= ((" + $MedianRef + " / 95.22216) * 100) - 100"
    $medianCell.ClearFormats()
}

# Low frequency group (data rows 2:31, mean in E3, median in E10).
Add-IncreaseBlock 18 19 "E3" "E10"

# Medium frequency group (data rows 34:63, mean in E35, median in E42).
Add-IncreaseBlock 50 51 "E35" "E42"

# High frequency group (data rows 66:95, mean in E68, median in E74).
Add-IncreaseBlock 82 83 "E68" "E74"

# "All" section header + overall average of the three group increases.
$hdrMean = $ws.Range("D113")
$hdrMean.Value = "Mean increase"
$hdrMean.Font.Bold = $true

$hdrMedian = $ws.Range("F113")
$hdrMedian.Value = "Median increase"
$hdrMedian.Font.Bold = $true

$meanAvg = $ws.Range("D114")
$meanAvg.Formula = "= (D19 + D51 + D83) / 3"
$meanAvg.ClearFormats()

$medianAvg = $ws.Range("F114")
$medianAvg.Formula = "= (F19 + F51 + F83) / 3"
$medianAvg.ClearFormats()

# Match the saved selection / scroll position from the authored edit.
$ws.Range("E98").Select()
